# need_to_buy.xlsx: clear the stale forecast figures for TISG (B), fcs (C),
# MYDIR (E) and need_to_buy (F) for every data row (2-15), leaving the
# "Giorno" (A) and "buy" (D) columns untouched - matching the refreshed
# R export which left those columns blank pending recomputation.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2:C15").ClearContents()
$ws.Range("B2:C15").Style = "Normal"

$ws.Range("E2:F15").ClearContents()
$ws.Range("E2:F15").Style = "Normal"
